# Inner JOIN appears to be working now
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: append two new rows (5 & 6) to the join-result table ---
$ws1.Range("A5").Value = 6
$ws1.Range("B5").Value = $true
$ws1.Range("C5").Value = '"Dad"'
$ws1.Range("D5").Value = 6

$ws1.Range("A6").Value = 7
$ws1.Range("B6").Value = $false
$ws1.Range("C6").Value = '"Mom"'
$ws1.Range("D6").Value = 6

# --- Sheet2: fix casing of "MOM" -> "Mom", add Column_3 of booleans, add a new data row ---
$ws2.Range("B3").Value = '"Mom"'

$ws2.Range("D1").Value = "Column_3"
$ws2.Range("D2").Value = $true
$ws2.Range("D3").Value = $true
$ws2.Range("D4").Value = $true
$ws2.Range("D5").Value = $true

$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = '"Mom"'
$ws2.Range("C6").Value = 8
$ws2.Range("D6").Value = $true

# --- Selections / active sheet state ---
$ws1.Range("D7").Select() | Out-Null
$ws2.Range("D2:D6").Select() | Out-Null

$ws2.Activate() | Out-Null
